$d = $word.ActiveDocument

# --- Paragraph 1: "Novo arquivo " + "word" (spell-checked) + "." -> single run "Novo arquivo word." ---
# Re-typing the same visible text via Find/Replace collapses the three runs
# (and drops the now-stale w:proofErr spell-check markers) into one run,
# exactly like Word does when it re-flows the paragraph after an edit.
$d.Content.Find.Execute("Novo arquivo word.", $true, $false, $false, $false, $false, $true, 1, $false, "Novo arquivo word.", 2)

# --- Paragraph 2: empty paragraph -> "Modificação rep1_1" ---
$modTxt = "Modifica" + [char]0x00E7 + [char]0x00E3 + "o rep1_1"
$p2 = $d.Paragraphs(2)
$p2.Range.InsertAfter($modTxt)
